$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "class name" row used to read "calss_name" (typo) - fix it to "example"
# before the row shifts happen below (it currently lives in row 10).
$ws.Cells.Item(10, 3).Value = "example"

# Insert a new row above row 8 to hold a "current filter" / "breadcrumb" locator
# (this pushes "collor filters" and everything below it down by one row).
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).Value = "current filter"
$ws.Cells.Item(8, 2).Value = "class name"
$ws.Cells.Item(8, 3).Value = "breadcrumb"

# Remove the now-obsolete "css selector" locator row (shifted down to row 12).
$ws.Rows.Item(12).Delete()

# Leave the selection where the author left off editing.
$ws.Range("C8").Select() | Out-Null
